# Update TPM-derived statistics for the Icam1-Il2ra LR-pair sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 9.861094666666666
$ws.Range("H2").Value = 29.583284
$ws.Range("I2").Value = 0.243709096397741
$ws.Range("J2").Value = 0.2437090963977409
$ws.Range("M2").Value = 0.3822983333333334
$ws.Range("N2").Value = 1.146895
$ws.Range("O2").Value = 0.04915201160487953
$ws.Range("P2").Value = 0.04915201160487953
$ws.Range("Q2").Value = 3.769880055908889
$ws.Range("R2").Value = 33.92892050318
$ws.Range("S2").Value = 0.01197879233435647
$ws.Range("T2").Value = 0.01197879233435647
$ws.Range("G3").Value = 9.861094666666666
$ws.Range("H3").Value = 29.583284
$ws.Range("I3").Value = 0.243709096397741
$ws.Range("J3").Value = 0.2437090963977409
$ws.Range("O3").Value = 0.3087451919724631
$ws.Range("P3").Value = 0.3087451919724631
$ws.Range("Q3").Value = 23.68025851986089
$ws.Range("R3").Value = 213.122326678748
$ws.Range("S3").Value = 0.07524401175275607
$ws.Range("T3").Value = 0.07524401175275605
$ws.Range("G4").Value = 9.861094666666666
$ws.Range("H4").Value = 29.583284
$ws.Range("I4").Value = 0.243709096397741
$ws.Range("J4").Value = 0.2437090963977409
$ws.Range("O4").Value = 0.6421027964226573
$ws.Range("P4").Value = 0.6421027964226573
$ws.Range("Q4").Value = 49.24824940098267
$ws.Range("R4").Value = 443.2342446088439
$ws.Range("S4").Value = 0.1564862923106284
$ws.Range("T4").Value = 0.1564862923106284
$ws.Range("I5").Value = 0.7254466225154019
$ws.Range("J5").Value = 0.7254466225154018
$ws.Range("M5").Value = 0.3822983333333334
$ws.Range("N5").Value = 1.146895
$ws.Range("O5").Value = 0.04915201160487953
$ws.Range("P5").Value = 0.04915201160487953
$ws.Range("Q5").Value = 11.22176723918389
$ws.Range("R5").Value = 100.995905152655
$ws.Range("S5").Value = 0.0356571608085977
$ws.Range("T5").Value = 0.03565716080859769
$ws.Range("I6").Value = 0.7254466225154019
$ws.Range("J6").Value = 0.7254466225154018
$ws.Range("O6").Value = 0.3087451919724631
$ws.Range("P6").Value = 0.3087451919724631
$ws.Range("S6").Value = 0.2239781567342928
$ws.Range("T6").Value = 0.2239781567342927
$ws.Range("I7").Value = 0.7254466225154019
$ws.Range("J7").Value = 0.7254466225154018
$ws.Range("O7").Value = 0.6421027964226573
$ws.Range("P7").Value = 0.6421027964226573
$ws.Range("S7").Value = 0.4658113049725114
$ws.Range("T7").Value = 0.4658113049725114
$ws.Range("I8").Value = 0.03084428108685718
$ws.Range("J8").Value = 0.03084428108685716
$ws.Range("M8").Value = 0.3822983333333334
$ws.Range("N8").Value = 1.146895
$ws.Range("O8").Value = 0.04915201160487953
$ws.Range("P8").Value = 0.04915201160487953
$ws.Range("Q8").Value = 0.4771231022022223
$ws.Range("R8").Value = 4.29410791982
$ws.Range("S8").Value = 0.00151605846192537
$ws.Range("T8").Value = 0.00151605846192537
$ws.Range("I9").Value = 0.03084428108685718
$ws.Range("J9").Value = 0.03084428108685716
$ws.Range("O9").Value = 0.3087451919724631
$ws.Range("P9").Value = 0.3087451919724631
$ws.Range("S9").Value = 0.009523023485414332
$ws.Range("T9").Value = 0.009523023485414329
$ws.Range("I10").Value = 0.03084428108685718
$ws.Range("J10").Value = 0.03084428108685716
$ws.Range("O10").Value = 0.6421027964226573
$ws.Range("P10").Value = 0.6421027964226573
$ws.Range("S10").Value = 0.01980519913951747
$ws.Range("T10").Value = 0.01980519913951747
